# Scheduled runner update: refresh currentAveragePrice / Leve profit figures
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets of the
# Chocobo_Profits workbook.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 28167536
$ws.Range("I132").Value = 35855124
$ws.Range("J132").Value = 1260975.4
$ws.Range("K132").Value = 107565372
$ws.Range("L132").Value = 3782926.2
$ws.Range("M132").Value = -107562842
$ws.Range("N132").Value = -3787986.2

$ws.Range("H137").Value = 955660.9399999999
$ws.Range("I137").Value = 1987452.4
$ws.Range("J137").Value = 3238.0386
$ws.Range("K137").Value = 5962357.199999999
$ws.Range("L137").Value = 9714.1158
$ws.Range("M137").Value = -5959807.199999999
$ws.Range("N137").Value = -14814.1158

$ws.Range("H138").Value = 3199.8
$ws.Range("I138").Value = 3199.8
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 9599.400000000001
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = -4459.400000000001
$ws.Range("N138").ClearContents()

$ws.Range("H140").Value = 63398.59
$ws.Range("J140").Value = 63398.59
$ws.Range("L140").Value = 63398.59
$ws.Range("N140").Value = -73758.59

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7111.625
$ws.Range("I32").Value = 9474.261
$ws.Range("J32").Value = 3915.1177
$ws.Range("K32").Value = 9474.261
$ws.Range("L32").Value = 3915.1177
$ws.Range("M32").Value = -9187.261
$ws.Range("N32").Value = -4489.1177

$ws.Range("H61").Value = 1521.6
$ws.Range("I61").Value = 1427
$ws.Range("K61").Value = 1427
$ws.Range("M61").Value = -1215

$ws.Range("H122").Value = 2949.75
$ws.Range("I122").Value = 1488.2
$ws.Range("J122").Value = 5385.6665
$ws.Range("K122").Value = 4464.6
$ws.Range("L122").Value = 16156.9995
$ws.Range("M122").Value = -2014.6
$ws.Range("N122").Value = -21056.9995

$ws.Range("H129").Value = 49867.6
$ws.Range("J129").Value = 49867.6
$ws.Range("L129").Value = 49867.6
$ws.Range("N129").Value = -59867.6

$ws.Range("H132").Value = 2756.139
$ws.Range("I132").Value = 1514.6818
$ws.Range("J132").Value = 4707
$ws.Range("K132").Value = 4544.0454
$ws.Range("L132").Value = 14121
$ws.Range("M132").Value = -2014.0454
$ws.Range("N132").Value = -19181

$ws.Range("H136").Value = 1521.6
$ws.Range("I136").Value = 1427
$ws.Range("K136").Value = 4281
$ws.Range("M136").Value = -1731

$ws.Range("H137").Value = 48367.6
$ws.Range("J137").Value = 48367.6
$ws.Range("L137").Value = 48367.6
$ws.Range("N137").Value = -58567.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3958.6
$ws.Range("I134").Value = 1624.6666
$ws.Range("J134").Value = 5868.1816
$ws.Range("K134").Value = 4873.9998
$ws.Range("L134").Value = 17604.5448
$ws.Range("M134").Value = -2338.9998
$ws.Range("N134").Value = -22674.5448

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 262290.53
$ws.Range("I31").Value = 644534.5
$ws.Range("J31").Value = 3351.0645
$ws.Range("K31").Value = 644534.5
$ws.Range("L31").Value = 3351.0645
$ws.Range("M31").Value = -644239.5
$ws.Range("N31").Value = -3941.0645

$ws.Range("H34").Value = 262290.53
$ws.Range("I34").Value = 644534.5
$ws.Range("J34").Value = 3351.0645
$ws.Range("K34").Value = 644534.5
$ws.Range("L34").Value = 3351.0645
$ws.Range("M34").Value = -644332.5
$ws.Range("N34").Value = -3755.0645

$ws.Range("H58").Value = 2444.7144
$ws.Range("I58").Value = 1350.8889
$ws.Range("K58").Value = 1350.8889
$ws.Range("M58").Value = -1147.8889

$ws.Range("H132").Value = 7191.2
$ws.Range("I132").Value = 6983.3335
$ws.Range("J132").Value = 7503
$ws.Range("K132").Value = 20950.0005
$ws.Range("L132").Value = 22509
$ws.Range("M132").Value = -18420.0005
$ws.Range("N132").Value = -27569

$ws.Range("H134").Value = 2247.625
$ws.Range("I134").Value = 997.25
$ws.Range("K134").Value = 2991.75
$ws.Range("M134").Value = -456.75

$ws.Range("H136").Value = 2444.7144
$ws.Range("I136").Value = 1350.8889
$ws.Range("K136").Value = 4052.6667
$ws.Range("M136").Value = -1502.6667

$ws.Range("H137").Value = 43580
$ws.Range("J137").Value = 43580
$ws.Range("L137").Value = 43580
$ws.Range("N137").Value = -53780

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 780.51
$ws.Range("I131").Value = 489.25
$ws.Range("J131").Value = 805.837
$ws.Range("K131").Value = 1467.75
$ws.Range("L131").Value = 2417.511
$ws.Range("M131").Value = 3572.25
$ws.Range("N131").Value = -12497.511

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5776.6523
$ws.Range("I122").Value = 2872
$ws.Range("J122").Value = 8945.362999999999
$ws.Range("K122").Value = 8616
$ws.Range("L122").Value = 26836.089
$ws.Range("M122").Value = -6166
$ws.Range("N122").Value = -31736.089

$ws.Range("H132").Value = 3585.1843
$ws.Range("I132").Value = 2522.3333
$ws.Range("J132").Value = 4898.1177
$ws.Range("K132").Value = 7566.999899999999
$ws.Range("L132").Value = 14694.3531
$ws.Range("M132").Value = -5036.999899999999
$ws.Range("N132").Value = -19754.3531

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1743.4286
$ws.Range("I61").Value = 1774.75
$ws.Range("J61").Value = 1701.6666
$ws.Range("K61").Value = 1774.75
$ws.Range("L61").Value = 1701.6666
$ws.Range("M61").Value = -1572.75
$ws.Range("N61").Value = -2105.6666

$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

$ws.Range("H113").Value = 1743.4286
$ws.Range("I113").Value = 1774.75
$ws.Range("J113").Value = 1701.6666
$ws.Range("K113").Value = 1774.75
$ws.Range("L113").Value = 1701.6666
$ws.Range("M113").Value = 395.25
$ws.Range("N113").Value = -6041.6666

$ws.Range("H132").Value = 8054.727
$ws.Range("I132").Value = 7268
$ws.Range("J132").Value = 8349.75
$ws.Range("K132").Value = 21804
$ws.Range("L132").Value = 25049.25
$ws.Range("M132").Value = -19274
$ws.Range("N132").Value = -30109.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 8784.333000000001
$ws.Range("I113").Value = 33700
$ws.Range("J113").Value = 479.1111
$ws.Range("K113").Value = 101100
$ws.Range("L113").Value = 1437.3333
$ws.Range("M113").Value = -98930
$ws.Range("N113").Value = -5777.3333

$ws.Range("H132").Value = 6176737
$ws.Range("I132").Value = 4661
$ws.Range("J132").Value = 12348813
$ws.Range("K132").Value = 13983
$ws.Range("L132").Value = 37046439
$ws.Range("M132").Value = -11453
$ws.Range("N132").Value = -37051499

Write-Host "Chocobo_Profits sheets updated"
